$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The greenish "git projects / quiz" callout box (cNvPr id="3", name
# "Rectangle 2") needs to move up and grow taller so the quiz bullet
# fits without overflowing.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 3) {
        $shp = $candidate
        break
    }
}

$shp.Top = 354.0000787401575
$shp.Height = 162
